$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in LoadTime_s values for the rows that were missing them
$ws.Range("H2").Value = 2
$ws.Range("H3").Value = 2
$ws.Range("H4").Value = 2

# Update the active selection to H5
$ws.Range("H5").Select()
